# AssumptionUsedByGroupDecision<N>  ->  AssumptionUsedDecision<N>
#
# Column A (rows 7..202) of the "Tariff" sheet holds 196 shared-string
# labels "AssumptionUsedByGroupDecision1" .. "AssumptionUsedByGroupDecision196".
# Rename them (in place) to "AssumptionUsedDecision1" .. "AssumptionUsedDecision196".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldPrefix = "AssumptionUsedByGroupDecision"
$newPrefix = "AssumptionUsedDecision"

for ($row = 7; $row -le 202; $row++) {
    $cell = $ws.Cells.Item($row, 1)
    $value = $cell.Value2
    if ($value -ne $null -and $value.ToString().StartsWith($oldPrefix)) {
        $suffix = $value.ToString().Substring($oldPrefix.Length)
        $cell.Value2 = $newPrefix + $suffix
    }
}

# Reflect the author's follow-up navigation: the active cell moved on to
# the next (empty) row right after the renamed list, A204.
$ws.Range("A204").Select() | Out-Null
